$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "2023-03-02-2013_RF_class_weights.csv"
$ws.Range("A16").Value = "2023-03-02-2054_RF_class_weights.csv"

$ws.Range("G15").Value = "0.620 (0.015)"
$ws.Range("G16").Value = "0.622 (0.010)"

$ws.Range("D15").Value = "max_features=1, n_jobs=6, class_weight={0:.35, 1:.65}"
$ws.Range("D16").Value = "max_features=1, n_jobs=6, class_weight={0:.1, 1:.9}"

$ws.Range("H15").Value = "March 2, 2023, 7:14 p.m."
$ws.Range("H16").Value = "March 2, 2023, 7:55 p.m."

$ws.Range("B15").Value = "RandomForest"
$ws.Range("B16").Value = "RandomForest"

$ws.Range("C15").Value = "MoCov"
$ws.Range("C16").Value = "MoCov"

$ws.Range("E15").Value = "1 x 3"
$ws.Range("E16").Value = "1 x 3"

$ws.Range("F15").Value = "weakly supervision with cv centers"
$ws.Range("F16").Value = "weakly supervision with cv centers"

$ws.Range("I15").Value = 0.598
$ws.Range("I16").Value = 0.594

$tbl = $ws.ListObjects.Item("Tabelle1")
$tbl.Resize($ws.Range("A1:I16"))

$ws.Range("H17").Select()
